$wb = $excel.ActiveWorkbook

# Work on the "Repayment Schedule" sheet
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (shifts N->O, O->P, P->Q)
$ws.Columns("N:N").Insert()

# Activate this sheet and select cell R5, making it the active tab
$ws.Activate()
$ws.Range("R5").Select()
